$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on the cells we are about to rewrite so that
# numeric-looking strings (e.g. "1.000", "233.05") are stored as text,
# matching the inlineStr cells in the source workbook.
$editRange = $ws.Range("D2:E51")
$editRange.NumberFormat = "@"

$ws.Range("D2").Value = "30.430.68"
$ws.Range("E2").Value = "  +0.94%  "
$ws.Range("D3").Value = "1.851.43"
$ws.Range("E3").Value = "  +1.00%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "233.05"
$ws.Range("E5").Value = "  +0.82%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").Value = "0.4739"
$ws.Range("E7").Value = "  +1.65%  "
$ws.Range("D8").Value = "0.2751"
$ws.Range("E8").Value = "  +2.49%  "
$ws.Range("D9").Value = "0.06341"
$ws.Range("E9").Value = "  +1.07%  "
$ws.Range("D10").Value = "17.90"
$ws.Range("E10").Value = "  +11.42%  "
$ws.Range("D11").Value = "1.860.77"
$ws.Range("E11").Value = "  +1.58%  "
$ws.Range("D12").Value = "0.07468"
$ws.Range("E12").Value = "  +1.10%  "
$ws.Range("D13").Value = "4.977"
$ws.Range("E13").Value = "  +1.62%  "
$ws.Range("D14").Value = "84.69"
$ws.Range("E14").Value = "  +1.61%  "
$ws.Range("D15").Value = "0.6230"
$ws.Range("E15").Value = "  +0.69%  "
$ws.Range("D16").Value = "30.380.11"
$ws.Range("E16").Value = "  +0.93%  "
$ws.Range("D17").Value = "246.52"
$ws.Range("E17").Value = "  +8.80%  "
$ws.Range("D18").Value = "1.000"
$ws.Range("E18").Value = "  -0.02%  "
$ws.Range("D19").Value = "12.68"
$ws.Range("E19").Value = "  +2.51%  "
$ws.Range("D20").Value = "0.000007335"
$ws.Range("E20").Value = "  +0.95%  "
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  +0.24%  "
$ws.Range("D22").Value = "4.930"
$ws.Range("D23").Value = "5.895"
$ws.Range("E23").Value = "  +0.62%  "
$ws.Range("D24").Value = "164.00"
$ws.Range("E24").Value = "  -0.41%  "
$ws.Range("D25").Value = "9.009"
$ws.Range("E25").Value = "  -1.74%  "
$ws.Range("D26").Value = "17.97"
$ws.Range("E26").Value = "  +2.10%  "
$ws.Range("E27").Value = "  +0.66%  "
$ws.Range("E28").Value = "  +0.98%  "
$ws.Range("D29").Value = "1.346"
$ws.Range("E29").Value = "  -1.55%  "
$ws.Range("D30").Value = "4.030"
$ws.Range("E30").Value = "  -0.64%  "
$ws.Range("D31").Value = "3.829"
$ws.Range("E31").Value = "  +1.37%  "
$ws.Range("D32").Value = "0.04818"
$ws.Range("E32").Value = "  +0.75%  "
$ws.Range("D33").Value = "1.127"
$ws.Range("E33").Value = "  -0.43%  "
$ws.Range("D34").Value = "0.6964"
$ws.Range("E34").Value = "  -1.31%  "
$ws.Range("D35").Value = "2.699"
$ws.Range("E35").Value = "  +0.59%  "
$ws.Range("D36").Value = "0.01892"
$ws.Range("E36").Value = "  +4.02%  "
$ws.Range("D37").Value = "2.677"
$ws.Range("E37").Value = "  +2.51%  "
$ws.Range("D38").Value = "0.8734"
$ws.Range("E38").Value = "  -2.25%  "
$ws.Range("D39").Value = "1.978"
$ws.Range("E39").Value = "  +2.40%  "
$ws.Range("D40").Value = "106.15"
$ws.Range("E40").Value = "  +2.57%  "
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("D42").Value = "0.4056"
$ws.Range("E42").Value = "  +1.41%  "
$ws.Range("D43").Value = "5.496"
$ws.Range("E43").Value = "  +0.68%  "
$ws.Range("D44").Value = "7.166"
$ws.Range("E44").Value = "  +2.71%  "
$ws.Range("D45").Value = "63.19"
$ws.Range("E45").Value = "  +5.87%  "
$ws.Range("D46").Value = "0.1197"
$ws.Range("E46").Value = "  +0.29%  "
$ws.Range("E47").Value = "  +4.02%  "
$ws.Range("D48").Value = "8.511"
$ws.Range("E48").Value = "  +0.36%  "
$ws.Range("D49").Value = "0.05498"
$ws.Range("E49").Value = "  -0.31%  "
$ws.Range("D50").Value = "1.347"
$ws.Range("E50").Value = "  -1.65%  "
$ws.Range("D51").Value = "0.3678"
$ws.Range("E51").Value = "  +1.45%  "

# Revert the temporary text format so the cells end up with the same
# (default) style they started with -- only their text content changed.
$editRange.ClearFormats()
